$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F-column "想去人数" (want-to-go count) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4
$ws1.Range("F5").Value = 495
$ws1.Range("F6").Value = 931
$ws1.Range("F7").Value = 167
$ws1.Range("F8").Value = 962
$ws1.Range("F9").Value = 752
$ws1.Range("F10").Value = 202
$ws1.Range("F15").Value = 557
$ws1.Range("F16").Value = 492
$ws1.Range("F17").Value = 1305
$ws1.Range("F19").Value = 434
$ws1.Range("F20").Value = 1121
$ws1.Range("F21").Value = 2809
$ws1.Range("F22").Value = 1320
$ws1.Range("F23").Value = 660
$ws1.Range("F24").Value = 170
$ws1.Range("F25").Value = 1249
$ws1.Range("F27").Value = 977
$ws1.Range("F28").Value = 323
$ws1.Range("F29").Value = 1495
$ws1.Range("F30").Value = 34
$ws1.Range("F32").Value = 1344

# Sheet "全部类型" (sheet4): F-column "想去人数" (want-to-go count) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4
$ws4.Range("F7").Value = 495
$ws4.Range("F13").Value = 931
$ws4.Range("F14").Value = 167
$ws4.Range("F16").Value = 962
$ws4.Range("F17").Value = 752
$ws4.Range("F18").Value = 202
$ws4.Range("F28").Value = 557
$ws4.Range("F29").Value = 492
$ws4.Range("F30").Value = 1305
$ws4.Range("F32").Value = 434
$ws4.Range("F33").Value = 1121
$ws4.Range("F34").Value = 2809
$ws4.Range("F35").Value = 1320
$ws4.Range("F36").Value = 660
$ws4.Range("F37").Value = 170
$ws4.Range("F38").Value = 1249
$ws4.Range("F42").Value = 977
$ws4.Range("F43").Value = 323
$ws4.Range("F44").Value = 1496
$ws4.Range("F45").Value = 34
$ws4.Range("F47").Value = 1344
